# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect the refreshed scrape data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 661
$wsExhibit.Range("F4").Value = 1525
$wsExhibit.Range("F5").Value = 705
$wsExhibit.Range("F6").Value = 16

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 661
$wsAll.Range("F4").Value = 1525
$wsAll.Range("F6").Value = 705
$wsAll.Range("F7").Value = 16
